$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1108.3334
$ws.Range("I33").Value = 1084.5555
$ws.Range("K33").Value = 1084.5555
$ws.Range("M33").Value = -855.5554999999999
$ws.Range("H62").Value = 4971.875
$ws.Range("I62").Value = 5196
$ws.Range("J62").Value = 4598.3335
$ws.Range("K62").Value = 5196
$ws.Range("L62").Value = 4598.3335
$ws.Range("M62").Value = -4572
$ws.Range("N62").Value = -5846.3335
$ws.Range("H65").Value = 4971.875
$ws.Range("I65").Value = 5196
$ws.Range("J65").Value = 4598.3335
$ws.Range("K65").Value = 25980
$ws.Range("L65").Value = 22991.6675
$ws.Range("M65").Value = -22860
$ws.Range("N65").Value = -29231.6675
$ws.Range("H74").Value = 10333.333
$ws.Range("I74").Value = 10000
$ws.Range("J74").Value = 10444.444
$ws.Range("K74").Value = 10000
$ws.Range("L74").Value = 10444.444
$ws.Range("M74").Value = -9064
$ws.Range("N74").Value = -12316.444
$ws.Range("H77").Value = 10333.333
$ws.Range("I77").Value = 10000
$ws.Range("J77").Value = 10444.444
$ws.Range("K77").Value = 50000
$ws.Range("L77").Value = 52222.22
$ws.Range("M77").Value = -45320
$ws.Range("N77").Value = -61582.22
$ws.Range("H80").Value = 3050.3462
$ws.Range("I80").Value = 1936
$ws.Range("J80").Value = 3460.8948
$ws.Range("K80").Value = 5808
$ws.Range("L80").Value = 10382.6844
$ws.Range("M80").Value = -4810
$ws.Range("N80").Value = -12378.6844
$ws.Range("H83").Value = 3050.3462
$ws.Range("I83").Value = 1936
$ws.Range("J83").Value = 3460.8948
$ws.Range("K83").Value = 17424
$ws.Range("L83").Value = 31148.0532
$ws.Range("M83").Value = -12432
$ws.Range("N83").Value = -41132.0532
$ws.Range("H88").Value = 3202
$ws.Range("I88").Value = 7048.2
$ws.Range("J88").Value = 1065.2222
$ws.Range("K88").Value = 7048.2
$ws.Range("L88").Value = 1065.2222
$ws.Range("M88").Value = -6642.2
$ws.Range("N88").Value = -1877.2222
$ws.Range("H91").Value = 3202
$ws.Range("I91").Value = 7048.2
$ws.Range("J91").Value = 1065.2222
$ws.Range("K91").Value = 7048.2
$ws.Range("L91").Value = 1065.2222
$ws.Range("M91").Value = -5644.2
$ws.Range("N91").Value = -3873.2222
$ws.Range("H100").Value = 5463.222
$ws.Range("I100").Value = 2353.75
$ws.Range("K100").Value = 2353.75
$ws.Range("M100").Value = -1812.75
$ws.Range("H101").Value = 1172.8462
$ws.Range("I101").Value = 460.77777
$ws.Range("J101").Value = 2775
$ws.Range("K101").Value = 1382.33331
$ws.Range("L101").Value = 8325
$ws.Range("M101").Value = 239.66669
$ws.Range("N101").Value = -11569
$ws.Range("H107").Value = 356.13333
$ws.Range("J107").Value = 363.14285
$ws.Range("L107").Value = 363.14285
$ws.Range("N107").Value = -4203.14285
$ws.Range("H116").Value = 3400
$ws.Range("I116").Value = 3400
$ws.Range("K116").Value = 3400
$ws.Range("M116").Value = 42
$ws.Range("H120").Value = 62999.668
$ws.Range("J120").Value = 62999.668
$ws.Range("L120").Value = 62999.668
$ws.Range("N120").Value = -72675.66800000001
$ws.Range("H127").Value = 8299
$ws.Range("I127").Value = 1881.6666
$ws.Range("K127").Value = 5644.9998
$ws.Range("M127").Value = -684.9997999999996
$ws.Range("H135").Value = 7578.4
$ws.Range("J135").Value = 18875
$ws.Range("L135").Value = 169875
$ws.Range("N135").Value = -174945

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9471.763000000001
$ws.Range("I32").Value = 7103.851
$ws.Range("K32").Value = 7103.851
$ws.Range("M32").Value = -6816.851
$ws.Range("H74").Value = 2073.3635
$ws.Range("I74").Value = 1822.1923
$ws.Range("K74").Value = 1822.1923
$ws.Range("M74").Value = -948.1922999999999
$ws.Range("H77").Value = 2073.3635
$ws.Range("I77").Value = 1822.1923
$ws.Range("K77").Value = 9110.961499999999
$ws.Range("M77").Value = -4742.961499999999
$ws.Range("H110").Value = 1939.85
$ws.Range("I110").Value = 2187.5625
$ws.Range("K110").Value = 2187.5625
$ws.Range("M110").Value = -142.5625
$ws.Range("H122").Value = 5524.9
$ws.Range("I122").Value = 4614.3335
$ws.Range("J122").Value = 6131.9443
$ws.Range("K122").Value = 13843.0005
$ws.Range("L122").Value = 18395.8329
$ws.Range("M122").Value = -11393.0005
$ws.Range("N122").Value = -23295.8329

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2012.56
$ws.Range("I94").Value = 1943.5238
$ws.Range("K94").Value = 1943.5238
$ws.Range("M94").Value = -1492.5238

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1169.5834
$ws.Range("I16").Value = 255.33333
$ws.Range("K16").Value = 255.33333
$ws.Range("M16").Value = 31.66667000000001
$ws.Range("H59").Value = 18540
$ws.Range("J59").Value = 22080
$ws.Range("L59").Value = 22080
$ws.Range("N59").Value = -24370
$ws.Range("H107").Value = 4822.5415
$ws.Range("I107").Value = 840.36365
$ws.Range("K107").Value = 840.36365
$ws.Range("M107").Value = 1079.63635
$ws.Range("H113").Value = 1169.5834
$ws.Range("I113").Value = 255.33333
$ws.Range("K113").Value = 255.33333
$ws.Range("M113").Value = 1914.66667
$ws.Range("H132").Value = 3227.4707
$ws.Range("I132").Value = 2451.6924
$ws.Range("K132").Value = 7355.0772
$ws.Range("M132").Value = -4825.0772
$ws.Range("H134").Value = 3585.0173
$ws.Range("I134").Value = 2527.0256
$ws.Range("K134").Value = 7581.0768
$ws.Range("M134").Value = -5046.0768

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 138526.27
$ws.Range("I63").Value = 401294
$ws.Range("J63").Value = 7142.4
$ws.Range("K63").Value = 1203882
$ws.Range("L63").Value = 21427.2
$ws.Range("M63").Value = -1203133
$ws.Range("N63").Value = -22925.2
$ws.Range("H66").Value = 138526.27
$ws.Range("I66").Value = 401294
$ws.Range("J66").Value = 7142.4
$ws.Range("K66").Value = 3611646
$ws.Range("L66").Value = 64281.6
$ws.Range("M66").Value = -3607902
$ws.Range("N66").Value = -71769.60000000001
$ws.Range("H98").Value = 1657.8
$ws.Range("I98").Value = 289
$ws.Range("K98").Value = 867
$ws.Range("M98").Value = 631
$ws.Range("H114").Value = 639.05
$ws.Range("J114").Value = 986.4
$ws.Range("L114").Value = 2959.2
$ws.Range("N114").Value = -9467.200000000001
$ws.Range("H139").Value = 7095
$ws.Range("J139").Value = 8990
$ws.Range("L139").Value = 26970
$ws.Range("N139").Value = -37250

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5039.857
$ws.Range("I122").Value = 2093
$ws.Range("K122").Value = 6279
$ws.Range("M122").Value = -3829

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2182.7585
$ws.Range("I61").Value = 1598.5714
$ws.Range("K61").Value = 1598.5714
$ws.Range("M61").Value = -1396.5714
$ws.Range("H109").Value = 48846.152
$ws.Range("J109").Value = 48846.152
$ws.Range("L109").Value = 48846.152
$ws.Range("N109").Value = -51620.152
$ws.Range("H113").Value = 2182.7585
$ws.Range("I113").Value = 1598.5714
$ws.Range("K113").Value = 1598.5714
$ws.Range("M113").Value = 571.4286

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 514.5
$ws.Range("I113").Value = 544
$ws.Range("K113").Value = 1632
$ws.Range("M113").Value = 538
